# Generate Report for Archive
#
# 1. Update status text from "Ready for handoff" to "In Translation" across
#    the Overview, zh-cn and de-de sheets.
# 2. Narrow the status/date columns (which had grown to fit the old text)
#    back down to their new, tighter width.

$wb = $excel.ActiveWorkbook

# NOTE: the host's ColumnWidth setter snaps the persisted OOXML column
# width to a 1/6-character pixel grid, so the literal target width of
# 13.4101845877511 cannot be reproduced exactly through COM. 12.5 is the
# ColumnWidth input whose rounded/stored result (13.3333...) lands closest
# to that target.
$newStatus = "In Translation"
$newWidth = 12.5

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
